$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---- Row: Piano | 31-32 | Only one set of dynamics needed for both staves,
#      since the material is in near-rhythmic unison, and doubled in octaves? | (empty)
$null = $t.Rows.Add()
$r = $t.Rows.Count
$t.Cell($r, 1).Range.Text = "Piano"
$t.Cell($r, 2).Range.Text = "31-32"
$t.Cell($r, 3).Range.Paragraphs.Item(1).Alignment = 0
$t.Cell($r, 3).Range.Text = "Only one set of dynamics needed for both staves, since the material is in near-rhythmic unison, and doubled in octaves?"

# ---- Row: Piano | 32 | Quaver stems should go down(?) since they are in bottom voice | (empty)
$null = $t.Rows.Add()
$r = $t.Rows.Count
$t.Cell($r, 1).Range.Text = "Piano"
$t.Cell($r, 2).Range.Text = "32"
$t.Cell($r, 3).Range.Paragraphs.Item(1).Alignment = 0
$t.Cell($r, 3).Range.Text = "Quaver stems should go down(?) since they are in bottom voice"

# ---- Row: Piano | 45 | Dynamic level (not in part) marked as pianissimo; inferred from
#      reference to other parts; clarify with composer | (empty)
$null = $t.Rows.Add()
$r = $t.Rows.Count
$t.Cell($r, 1).Range.Text = "Piano"
$t.Cell($r, 2).Range.Text = "45"
$t.Cell($r, 3).Range.Paragraphs.Item(1).Alignment = 0
$t.Cell($r, 3).Range.Text = "Dynamic level (not in part) marked as pianissimo; inferred from reference to other parts; clarify with composer"
